$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ------------------------------------------------------------------
# Sheet "汽車" (car) currently only has columns B:G and its header row
# is a stray copy of row 2's data. This change turns row 1 into a real
# header, fixes the two data rows, and extends every row out to column N
# with the same legislator/source metadata columns used on the other
# property sheets.
# ------------------------------------------------------------------

# --- Row 1: real header labels for B1:G1, new headers H1:N1 ---
# Previously B1:G1 were a stray duplicate of row 2's data values; replace
# them with the actual column names.
$ws.Range("B1").Value2 = "name"
$ws.Range("C1").Value2 = "capacity"
$ws.Range("D1").Value2 = "owner"
$ws.Range("E1").Value2 = "register_date"
$ws.Range("F1").Value2 = "register_reason"
$ws.Range("G1").Value2 = "acquire_value"

$ws.Range("B1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H1").Value2 = "property_category"
$ws.Range("I1").Value2 = "category"
$ws.Range("J1").Value2 = "date"
$ws.Range("K1").Value2 = "legislator_name"
$ws.Range("L1").Value2 = "legislator_id"
$ws.Range("M1").Value2 = "source_file"
$ws.Range("N1").Value2 = "index"

# --- Row 2: TOYOTA car (unchanged core data), extended with metadata ---
$ws.Range("B2").Value2 = "TOYOTA(汽車）"
$ws.Range("C2").Value2 = 2995
$ws.Range("E2").Value2 = "93年08月05日"
$ws.Range("G2").Value2 = 1090000

$ws.Range("H2").NumberFormat = "@"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("H2").Value2 = "land"
$ws.Range("I2").Value2 = "normal"
$ws.Range("J2").Value2 = "2013-07-11"
$ws.Range("K2").Value2 = "林正二"
$ws.Range("L2").Value2 = 788
$ws.Range("M2").Value2 = "tmp685a1"
$ws.Range("N2").Value2 = 40

$ws.Range("B2").Copy()
$ws.Range("H2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 3: NISSAN car, extended with metadata ---
$ws.Range("B3").Value2 = "NISSAN(汽車）"
$ws.Range("E3").Value2 = "99年02月09日"
$ws.Range("G3").Value2 = 1000000

$ws.Range("H3").NumberFormat = "@"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("H3").Value2 = "land"
$ws.Range("I3").Value2 = "normal"
$ws.Range("J3").Value2 = "2013-07-11"
$ws.Range("K3").Value2 = "林正二"
$ws.Range("L3").Value2 = 788
$ws.Range("M3").Value2 = "tmp685a1"
$ws.Range("N3").Value2 = 41

$ws.Range("B3").Copy()
$ws.Range("H3:N3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
